$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "1.0389 at -121.85"
$ws.Range("D3").Value = "1.0140 at 117.90"

$ws.Range("C4").Value = "1.0372 at -121.92"
$ws.Range("D4").Value = "1.0119 at 117.95"

$ws.Range("B5").Value = "1.0209 at -2.49"
$ws.Range("C5").Value = "1.0419 at -121.72"
$ws.Range("D5").Value = "1.0173 at 117.83"

$ws.Range("B6").Value = "1.0202 at -2.51"
$ws.Range("C6").Value = "1.0414 at -121.73"
$ws.Range("D6").Value = "1.0168 at 117.82"

$ws.Range("B7").Value = "1.0044 at -3.48"
$ws.Range("C7").Value = "1.0423 at -122.75"
$ws.Range("D7").Value = "1.0175 at 117.05"

$ws.Range("B8").Value = "1.0015 at -3.48"
$ws.Range("D8").Value = "1.0188 at 117.01"

$ws.Range("D9").Value = "1.0202 at 116.93"

$ws.Range("B10").Value = "1.0055 at -3.58"
$ws.Range("C10").Value = "1.0436 at -122.84"
$ws.Range("D10").Value = "1.0190 at 116.97"

$ws.Range("B11").Value = "1.0044 at -3.48"
$ws.Range("C11").Value = "1.0423 at -122.75"
$ws.Range("D11").Value = "1.0175 at 117.05"

$ws.Range("B12").Value = "0.9959 at -3.41"
